$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Azerbaijan Premier League")

# Row 149 (id 147): 2024-04-14 09:30 Sabail FC vs FK Kapaz, FTR=D
$row149 = @(
    147, 7011640, "Azerbaijan Premier League", "Azerbaijan Premier League", 45396.39583333334,
    "Sabail FC", "FK Kapaz", 3, 3, "D",
    1.727, 3.25, 4.333, 1.7, 3.8, 3.8,
    -0.75, 1.925, 1.875, 3.25, 1.775, 1.925,
    -1, 2.8, -1, -1, 0.875, 0.7749999999999999, -1
)

# Row 150 (id 148): 2024-04-14 12:00 Sabah vs FK Qarabag, FTR=H
$row150 = @(
    148, 7012356, "Azerbaijan Premier League", "Azerbaijan Premier League", 45396.5,
    "Sabah", "FK Qarabag", 3, 2, "H",
    2.625, 4, 2.1, 4, 4, 1.615,
    0.75, 1.975, 1.825, 2.75, 1.85, 1.95,
    3, -1, -1, 0.9750000000000001, -1, 0.8500000000000001, -1
)

for ($i = 0; $i -lt $row149.Count; $i++) {
    $ws.Cells.Item(149, $i + 1).Value = $row149[$i]
}

for ($i = 0; $i -lt $row150.Count; $i++) {
    $ws.Cells.Item(150, $i + 1).Value = $row150[$i]
}

# Apply formatting matching the rest of the table:
# Column A (id) uses bold font + thin border + centered/top alignment.
# Column E (Date) uses the custom "YYYY-MM-DD HH:MM:SS" number format.
foreach ($r in 149, 150) {
    $idCell = $ws.Cells.Item($r, 1)
    $idCell.Font.Bold = $true
    $idCell.Borders.LineStyle = 1
    $idCell.HorizontalAlignment = -4108
    $idCell.VerticalAlignment = -4160

    $dateCell = $ws.Cells.Item($r, 5)
    $dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
